$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the entire row 642 ("恥ずかしがり屋のマヌルネコと暖かい朝の太陽光線" post),
# shifting all subsequent rows up by one.
$ws.Rows.Item(642).Delete()
